$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card1")

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 15).Value = "nan"
}
